$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Snippets" table currently spans A1:D184. We are adding a new sample
# entry describing the worksheet autofilter snippet, which takes two new
# table rows (185 and 186).
$lo = $ws.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Fill in the new rows. Columns are entered A, B, D, C (rather than strict
# left-to-right order) so that newly introduced shared strings are created
# in the same order as the target workbook.
$ws.Range("A185").Value = "Worksheet"
$ws.Range("B185").Value = "autofilter"
$ws.Range("D185").Value = "addAutoFilter"
$ws.Range("C185").Value = "excel-worksheet-auto-filter"

$ws.Range("A186").Value = "AutoFilter"
$ws.Range("B186").Value = "apply"
$ws.Range("D186").Value = "addAutoFilter"
$ws.Range("C186").Value = "excel-worksheet-auto-filter"

# Row 184 was previously the most-recently-added sample and carried a
# highlight style; that highlight now belongs to the new rows, so clear it
# from row 184.
$ws.Range("A184:D184").ClearFormats()

# Move the active selection to reflect where the author ended up editing.
$ws.Range("O178").Select() | Out-Null
